$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.537.13'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.12%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.989.39'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.87%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9991'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4674'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.61%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3939'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.53%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07944'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.94%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.002'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.60%  '

$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.90'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.19%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.977.75'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.21%  '

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.248'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.62%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.867'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.93%  '

$ws.Range("B15").Value = 'TRON'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07127'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.61%  '

$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.61'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.73%  '

$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.08%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009949'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.46%  '

$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.61%  '

$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9990'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.24%  '

$ws.Range("B21").Value = 'WrappedBTC'
$ws.Range("C21").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.582.69'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.28%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.529'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.54%  '

$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.27'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.81%  '

$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.218.80'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.14%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.100'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.06%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.62'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.89%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.65'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.968'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.49%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '120.33'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.58%  '

$ws.Range("E30").Value = '  +1.86%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09444'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9057'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.29%  '

$ws.Range("E33").Value = '  +0.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.349'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.175'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.47%  '

$ws.Range("E36").Value = '  +1.47%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.174'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.49%  '

$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.000003375'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +106.30%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02119'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.46%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.909'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5751'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.78%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1823'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.35%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.815'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.37%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.03'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.40%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5371'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.55%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.680'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.97%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.186'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.18%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06946'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.42%  '

$ws.Range("E49").Value = '  +1.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '114.39'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3091'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.05%  '
